$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format numeric-looking Price cells as Text so Excel keeps them as strings
$textCells = @("D5", "D6", "D10", "D20", "D22", "D24", "D28", "D29", "D32", "D34", "D35", "D40", "D41", "D42", "D45", "D46", "D47", "D50", "D51")
foreach ($cellAddr in $textCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "66.849.28"
$ws.Range("E2").Value = "  -2.89%  "
$ws.Range("D3").Value = "3.472.06"
$ws.Range("E3").Value = "  -2.76%  "
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").Value = "604.75"
$ws.Range("E5").Value = "  -2.74%  "
$ws.Range("D6").Value = "148.60"
$ws.Range("E6").Value = "  -5.02%  "
$ws.Range("D7").Value = "3.470.48"
$ws.Range("E7").Value = "  -2.81%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("D10").Value = "0.143"
$ws.Range("E10").Value = "  -3.38%  "
$ws.Range("E11").Value = "  +3.04%  "
$ws.Range("E12").Value = "  -3.20%  "
$ws.Range("E13").Value = "  -4.11%  "
$ws.Range("E14").Value = "  -4.49%  "
$ws.Range("D15").Value = "4.059.43"
$ws.Range("E15").Value = "  -2.92%  "
$ws.Range("D16").Value = "3.473.42"
$ws.Range("E16").Value = "  -3.04%  "
$ws.Range("D17").Value = "66.818.77"
$ws.Range("E17").Value = "  -3.64%  "
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("E19").Value = "  -4.45%  "
$ws.Range("D20").Value = "15.45"
$ws.Range("E20").Value = "  -3.70%  "
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("D22").Value = "440.68"
$ws.Range("E22").Value = "  -4.17%  "
$ws.Range("E23").Value = "  -4.52%  "
$ws.Range("D24").Value = "79.37"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").Value = "3.608.50"
$ws.Range("E26").Value = "  -3.16%  "
$ws.Range("E27").Value = "  -8.79%  "
$ws.Range("D28").Value = "9.81"
$ws.Range("E28").Value = "  -7.54%  "
$ws.Range("D29").Value = "8.39"
$ws.Range("E29").Value = "  -7.98%  "
$ws.Range("E30").Value = "  -4.76%  "
$ws.Range("E31").Value = "  -6.07%  "
$ws.Range("D32").Value = "0.168"
$ws.Range("E32").Value = "  -2.28%  "
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").Value = "25.46"
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("D35").Value = "6.08"
$ws.Range("E35").Value = "  -5.90%  "
$ws.Range("D36").Value = "3.461.95"
$ws.Range("E36").Value = "  -3.05%  "
$ws.Range("E37").Value = "  -5.97%  "
$ws.Range("E38").Value = "  -4.74%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("D41").Value = "177.08"
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("D42").Value = "0.0891"
$ws.Range("E42").Value = "  -3.02%  "
$ws.Range("E43").Value = "  -10.23%  "
$ws.Range("E44").Value = "  -3.64%  "
$ws.Range("D45").Value = "0.888"
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("D46").Value = "29.27"
$ws.Range("E46").Value = "  -5.23%  "
$ws.Range("D47").Value = "46.25"
$ws.Range("E47").Value = "  +0.98%  "
$ws.Range("E48").Value = "  -9.25%  "
$ws.Range("E49").Value = "  -8.84%  "
$ws.Range("D50").Value = "7.48"
$ws.Range("E50").Value = "  -3.98%  "
$ws.Range("D51").Value = "0.990"
$ws.Range("E51").Value = "  -4.02%  "

# Restore default style on cells we temporarily formatted as Text
foreach ($cellAddr in $textCells) {
    $ws.Range($cellAddr).Style = "Normal"
}
